{"js": "// Update the worksheet date and every two-digit multiplication problem.\n// Every <w:t> run in the document changes, so we do a literal\n// find-and-replace for each (old -> new) pair. All old strings are\n// unique in the document, so `body.search()` unambiguously locates the\n// single run to rewrite.\nconst replacements = [\n  [\"2024-09-13 Friday\", \"2024-09-14 Saturday\"],\n  [\"79\u00d737=\", \"71\u00d728=\"],\n  [\"59\u00d774=\", \"47\u00d738=\"],\n  [\"48\u00d734=\", \"68\u00d786=\"],\n  [\"74\u00d714=\", \"73\u00d766=\"],\n  [\"20\u00d782=\", \"50\u00d747=\"],\n  [\"47\u00d739=\", \"33\u00d735=\"],\n  [\"85\u00d757=\", \"62\u00d738=\"],\n  [\"63\u00d736=\", \"63\u00d778=\"],\n  [\"87\u00d763=\", \"95\u00d760=\"],\n  [\"94\u00d796=\", \"59\u00d728=\"],\n  [\"50\u00d751=\", \"57\u00d743=\"],\n  [\"28\u00d723=\", \"96\u00d721=\"],\n  [\"70\u00d768=\", \"59\u00d716=\"],\n  [\"32\u00d713=\", \"89\u00d780=\"],\n  [\"37\u00d752=\", \"55\u00d777=\"],\n  [\"42\u00d740=\", \"82\u00d759=\"],\n  [\"85\u00d730=\", \"58\u00d763=\"],\n  [\"70\u00d738=\", \"38\u00d785=\"],\n  [\"44\u00d773=\", \"95\u00d759=\"],\n  [\"27\u00d718=\", \"96\u00d797=\"],\n  [\"94\u00d734=\", \"32\u00d751=\"],\n  [\"80\u00d771=\", \"15\u00d745=\"],\n  [\"29\u00d748=\", \"28\u00d780=\"],\n  [\"68\u00d773=\", \"77\u00d736=\"],\n  [\"86\u00d763=\", \"14\u00d722=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true });\n  found.load(\"text\");\n  await context.sync();\n\n  for (let i = 0; i < found.items.length; i++) {\n    found.items[i].insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the worksheet date and every two-digit multiplication problem.\n# Every <w:t> run in the document changes, so we run a literal\n# (non-wildcard) Find/Replace for each (old -> new) pair across the whole\n# document body. All old strings are unique, so wdReplaceAll only ever\n# touches a single run per pair.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-09-13 Friday\", \"2024-09-14 Saturday\"),\n    @(\"79\u00d737=\", \"71\u00d728=\"),\n    @(\"59\u00d774=\", \"47\u00d738=\"),\n    @(\"48\u00d734=\", \"68\u00d786=\"),\n    @(\"74\u00d714=\", \"73\u00d766=\"),\n    @(\"20\u00d782=\", \"50\u00d747=\"),\n    @(\"47\u00d739=\", \"33\u00d735=\"),\n    @(\"85\u00d757=\", \"62\u00d738=\"),\n    @(\"63\u00d736=\", \"63\u00d778=\"),\n    @(\"87\u00d763=\", \"95\u00d760=\"),\n    @(\"94\u00d796=\", \"59\u00d728=\"),\n    @(\"50\u00d751=\", \"57\u00d743=\"),\n    @(\"28\u00d723=\", \"96\u00d721=\"),\n    @(\"70\u00d768=\", \"59\u00d716=\"),\n    @(\"32\u00d713=\", \"89\u00d780=\"),\n    @(\"37\u00d752=\", \"55\u00d777=\"),\n    @(\"42\u00d740=\", \"82\u00d759=\"),\n    @(\"85\u00d730=\", \"58\u00d763=\"),\n    @(\"70\u00d738=\", \"38\u00d785=\"),\n    @(\"44\u00d773=\", \"95\u00d759=\"),\n    @(\"27\u00d718=\", \"96\u00d797=\"),\n    @(\"94\u00d734=\", \"32\u00d751=\"),\n    @(\"80\u00d771=\", \"15\u00d745=\"),\n    @(\"29\u00d748=\", \"28\u00d780=\"),\n    @(\"68\u00d773=\", \"77\u00d736=\"),\n    @(\"86\u00d763=\", \"14\u00d722=\")\n)\n\nforeach ($pair in $replacements) {\n    $old = $pair[0]\n    $new = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $old\n    $find.Replacement.Text = $new\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Forward = $true\n    $find.Wrap = 1\n\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null\n}\n"}
